$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1243.6842

$ws.Range("H93").Value = 21249.166
$ws.Range("J93").Value = 21249.166
$ws.Range("L93").Value = 21249.166
$ws.Range("N93").Value = -26241.166

$ws.Range("H98").Value = 9322.923000000001
$ws.Range("I98").Value = 6000
$ws.Range("J98").Value = 9927.091
$ws.Range("K98").Value = 6000
$ws.Range("L98").Value = 9927.091
$ws.Range("M98").Value = -4502
$ws.Range("N98").Value = -12923.091

$ws.Range("H107").Value = 2098.182
$ws.Range("I107").Value = 2098.182
$ws.Range("K107").Value = 2098.182
$ws.Range("M107").Value = -178.1819999999998

$ws.Range("H112").Value = 1279.3281
$ws.Range("J112").Value = 1279.3281
$ws.Range("L112").Value = 3837.9843
$ws.Range("N112").Value = -6053.9843

$ws.Range("H122").Value = 9322.923000000001
$ws.Range("I122").Value = 6000
$ws.Range("J122").Value = 9927.091
$ws.Range("K122").Value = 18000
$ws.Range("L122").Value = 29781.273
$ws.Range("M122").Value = -15550
$ws.Range("N122").Value = -34681.273

$ws.Range("H138").Value = 2560.2632
$ws.Range("I138").Value = 1399
$ws.Range("J138").Value = 3164.12
$ws.Range("K138").Value = 4197
$ws.Range("L138").Value = 9492.360000000001
$ws.Range("M138").Value = 943
$ws.Range("N138").Value = -19772.36

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1451.697
$ws.Range("I2").Value = 1458.6154
$ws.Range("J2").Value = 1426
$ws.Range("K2").Value = 1458.6154
$ws.Range("L2").Value = 1426
$ws.Range("M2").Value = -1345.6154
$ws.Range("N2").Value = -1652

$ws.Range("H35").Value = 18599.75
$ws.Range("J35").Value = 28799.5
$ws.Range("L35").Value = 28799.5
$ws.Range("N35").Value = -29611.5

$ws.Range("H36").Value = 14253.333
$ws.Range("J36").Value = 16750
$ws.Range("L36").Value = 16750
$ws.Range("N36").Value = -17442

$ws.Range("H115").Value = 29956.666
$ws.Range("J115").Value = 29956.666
$ws.Range("L115").Value = 29956.666
$ws.Range("N115").Value = -33090.666

$ws.Range("H116").Value = 1451.697
$ws.Range("I116").Value = 1458.6154
$ws.Range("J116").Value = 1426
$ws.Range("K116").Value = 1458.6154
$ws.Range("L116").Value = 1426
$ws.Range("M116").Value = 835.3846000000001
$ws.Range("N116").Value = -6014

$ws.Range("H118").Value = 28390
$ws.Range("J118").Value = 28390
$ws.Range("L118").Value = 28390
$ws.Range("N118").Value = -31704

$ws.Range("H119").Value = 35465.332
$ws.Range("J119").Value = 35465.332
$ws.Range("L119").Value = 35465.332
$ws.Range("N119").Value = -45141.332

$ws.Range("H132").Value = 2390.283
$ws.Range("I132").Value = 1472.3667
$ws.Range("J132").Value = 3587.5652
$ws.Range("K132").Value = 4417.1001
$ws.Range("L132").Value = 10762.6956
$ws.Range("M132").Value = -1887.1001
$ws.Range("N132").Value = -15822.6956

$ws.Range("H137").Value = 59999
$ws.Range("J137").Value = 59999
$ws.Range("L137").Value = 59999
$ws.Range("N137").Value = -70199

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1451.697
$ws.Range("I3").Value = 1458.6154
$ws.Range("J3").Value = 1426
$ws.Range("K3").Value = 1458.6154
$ws.Range("L3").Value = 1426
$ws.Range("M3").Value = -1344.6154
$ws.Range("N3").Value = -1654

$ws.Range("H36").Value = 20000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 20000
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -21068

$ws.Range("H95").Value = 33285.715
$ws.Range("J95").Value = 33285.715
$ws.Range("L95").Value = 33285.715
$ws.Range("N95").Value = -38777.715

$ws.Range("H99").Value = 1931.5385
$ws.Range("I99").Value = 1842.5
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 1842.5
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -344.5
$ws.Range("N99").Value = -5996

$ws.Range("H134").Value = 3826.3044
$ws.Range("I134").Value = 1554.5834
$ws.Range("K134").Value = 4663.7502
$ws.Range("M134").Value = -2128.7502

$ws.Range("H137").Value = 50660
$ws.Range("J137").Value = 50660
$ws.Range("L137").Value = 50660
$ws.Range("N137").Value = -60860

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2571.0908
$ws.Range("I58").Value = 1516.3182
$ws.Range("J58").Value = 4680.636
$ws.Range("K58").Value = 1516.3182
$ws.Range("L58").Value = 4680.636
$ws.Range("M58").Value = -1313.3182
$ws.Range("N58").Value = -5086.636

$ws.Range("H132").Value = 4631.9287
$ws.Range("I132").Value = 3834.3845
$ws.Range("J132").Value = 15000
$ws.Range("K132").Value = 11503.1535
$ws.Range("L132").Value = 45000
$ws.Range("M132").Value = -8973.1535
$ws.Range("N132").Value = -50060

$ws.Range("H136").Value = 2571.0908
$ws.Range("I136").Value = 1516.3182
$ws.Range("J136").Value = 4680.636
$ws.Range("K136").Value = 4548.9546
$ws.Range("L136").Value = 14041.908
$ws.Range("M136").Value = -1998.9546
$ws.Range("N136").Value = -19141.908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 82.3125
$ws.Range("I12").Value = 11
$ws.Range("J12").Value = 92.5
$ws.Range("K12").Value = 33
$ws.Range("L12").Value = 277.5
$ws.Range("M12").Value = 140
$ws.Range("N12").Value = -623.5

$ws.Range("H33").Value = 132.25
$ws.Range("J33").Value = 149.21428
$ws.Range("L33").Value = 895.28568
$ws.Range("N33").Value = -1461.28568

$ws.Range("H68").Value = 1142.1578
$ws.Range("I68").Value = 678.9259
$ws.Range("J68").Value = 2279.182
$ws.Range("K68").Value = 2036.7777
$ws.Range("L68").Value = 6837.545999999999
$ws.Range("M68").Value = -1225.7777
$ws.Range("N68").Value = -8459.545999999998

$ws.Range("H71").Value = 1142.1578
$ws.Range("I71").Value = 678.9259
$ws.Range("J71").Value = 2279.182
$ws.Range("K71").Value = 6110.3331
$ws.Range("L71").Value = 20512.638
$ws.Range("M71").Value = -2054.3331
$ws.Range("N71").Value = -28624.638

$ws.Range("H113").Value = 652.4838999999999
$ws.Range("I113").Value = 573.9
$ws.Range("J113").Value = 689.9048
$ws.Range("K113").Value = 1721.7
$ws.Range("L113").Value = 2069.7144
$ws.Range("M113").Value = 448.3000000000002
$ws.Range("N113").Value = -6409.7144

$ws.Range("H123").Value = 6990.909
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 6990.909
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 20972.727
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -25872.727

$ws.Range("H131").Value = 785.22
$ws.Range("J131").Value = 815.8022
$ws.Range("L131").Value = 2447.4066
$ws.Range("N131").Value = -12527.4066

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 25641.908
$ws.Range("J46").Value = 26706.6
$ws.Range("L46").Value = 26706.6
$ws.Range("N46").Value = -27018.6

$ws.Range("H70").Value = 6287.946
$ws.Range("I70").Value = 5922.6895
$ws.Range("J70").Value = 7612
$ws.Range("K70").Value = 5922.6895
$ws.Range("L70").Value = 7612
$ws.Range("M70").Value = -5652.6895
$ws.Range("N70").Value = -8152

$ws.Range("H73").Value = 6287.946
$ws.Range("I73").Value = 5922.6895
$ws.Range("J73").Value = 7612
$ws.Range("K73").Value = 5922.6895
$ws.Range("L73").Value = 7612
$ws.Range("M73").Value = -4986.6895
$ws.Range("N73").Value = -9484

$ws.Range("H113").Value = 1301.625
$ws.Range("I113").Value = 1231.1666
$ws.Range("J113").Value = 1513
$ws.Range("K113").Value = 1231.1666
$ws.Range("L113").Value = 1513
$ws.Range("M113").Value = 938.8334
$ws.Range("N113").Value = -5853

$ws.Range("H124").Value = 41880
$ws.Range("J124").Value = 41880
$ws.Range("L124").Value = 41880
$ws.Range("N124").Value = -51700

$ws.Range("H126").Value = 3916.3674
$ws.Range("I126").Value = 2911.4856
$ws.Range("K126").Value = 8734.4568
$ws.Range("M126").Value = -6264.4568

$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 3897
$ws.Range("I132").Value = 2868.0588
$ws.Range("J132").Value = 5354.6665
$ws.Range("K132").Value = 8604.1764
$ws.Range("L132").Value = 16063.9995
$ws.Range("M132").Value = -6074.1764
$ws.Range("N132").Value = -21123.9995

$ws.Range("H137").Value = 20000
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1972.5
$ws.Range("I61").Value = 1963.3334
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1963.3334
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1761.3334
$ws.Range("N61").Value = -2404

$ws.Range("H109").Value = 19250
$ws.Range("J109").Value = 19250
$ws.Range("L109").Value = 19250
$ws.Range("N109").Value = -22024

$ws.Range("H113").Value = 1972.5
$ws.Range("I113").Value = 1963.3334
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1963.3334
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 206.6666
$ws.Range("N113").Value = -6340

$ws.Range("H122").Value = 3874.261
$ws.Range("I122").Value = 2676
$ws.Range("J122").Value = 4126.5264
$ws.Range("K122").Value = 8028
$ws.Range("L122").Value = 12379.5792
$ws.Range("M122").Value = -5578
$ws.Range("N122").Value = -17279.5792

$ws.Range("H132").Value = 8100.8
$ws.Range("I132").Value = 10104
$ws.Range("J132").Value = 7600
$ws.Range("K132").Value = 30312
$ws.Range("L132").Value = 22800
$ws.Range("M132").Value = -27782
$ws.Range("N132").Value = -27860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 9812.5
$ws.Range("J74").Value = 9875
$ws.Range("L74").Value = 9875
$ws.Range("N74").Value = -11747

$ws.Range("H77").Value = 9812.5
$ws.Range("J77").Value = 9875
$ws.Range("L77").Value = 29625
$ws.Range("N77").Value = -38985
